$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'290.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.06%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'30.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-1.67%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.871"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.59%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07246"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-0.14%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.346"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'32.43%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.648"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.00%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.706"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.84%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.8974"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-1.25%"
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'2.80%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.08067"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'8.03%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08114"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.09%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03077"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'2.44%"
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'0.17%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001505"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.12%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005792"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.07%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.482"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.95%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.074"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-2.31%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3317"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.75%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1288"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.68%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'3.962"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-9.72%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'5.27%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04517"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.03%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-2.27%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004404"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'10.53%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001299"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'2.75%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D39").Value = "'0.01586"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-1.64%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04379"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.23%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007254"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-1.53%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.01000"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.1312"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'0.26%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002078"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-3.79%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009167"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-16.91%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00005722"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-6.35%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-1.33%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'25.76%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'-4.64%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-1.33%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-1.33%"
$ws.Range("E51").Style = "Normal"
